# Refresh cached market-board profit figures on the Asura_Profits workbook
# (scheduled Sheets runner: pulls latest Universalis prices and rewrites
#  the currentAveragePrice* / LevePrice* / LeveProfit* columns per leve row)

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 306.58334
$ws.Range("I2").Value = 297.66666
$ws.Range("J2").Value = 333.33334
$ws.Range("K2").Value = 297.66666
$ws.Range("L2").Value = 333.33334
$ws.Range("M2").Value = -184.66666
$ws.Range("N2").Value = -559.33334
# Row 62
$ws.Range("H62").Value = 4065.2307
$ws.Range("I62").Value = 2883.5557
$ws.Range("J62").Value = 6724
$ws.Range("K62").Value = 2883.5557
$ws.Range("L62").Value = 6724
$ws.Range("M62").Value = -2259.5557
$ws.Range("N62").Value = -7972
# Row 65
$ws.Range("H65").Value = 4065.2307
$ws.Range("I65").Value = 2883.5557
$ws.Range("J65").Value = 6724
$ws.Range("K65").Value = 14417.7785
$ws.Range("L65").Value = 33620
$ws.Range("M65").Value = -11297.7785
$ws.Range("N65").Value = -39860
# Row 98
$ws.Range("H98").Value = 4353.3335
$ws.Range("I98").Value = 4353.3335
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 4353.3335
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -2855.3335
# Row 113
$ws.Range("H113").Value = 3193.182
$ws.Range("I113").Value = 2401.25
$ws.Range("J113").Value = 3645.7144
$ws.Range("K113").Value = 2401.25
$ws.Range("L113").Value = 3645.7144
$ws.Range("M113").Value = 852.75
$ws.Range("N113").Value = -10153.7144
# Row 122
$ws.Range("H122").Value = 4353.3335
$ws.Range("I122").Value = 4353.3335
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13060.0005
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10610.0005
# Row 137
$ws.Range("H137").Value = 1265.4445
$ws.Range("I137").Value = 1221.8718
$ws.Range("J137").Value = 1548.6666
$ws.Range("K137").Value = 3665.6154
$ws.Range("L137").Value = 4645.9998
$ws.Range("M137").Value = -1115.6154
$ws.Range("N137").Value = -9745.9998

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12681.883
$ws.Range("I32").Value = 14203
$ws.Range("J32").Value = 7150.5454
$ws.Range("K32").Value = 14203
$ws.Range("L32").Value = 7150.5454
$ws.Range("M32").Value = -13916
# Row 61
$ws.Range("H61").Value = 3975.0908
$ws.Range("I61").Value = 5047.2
$ws.Range("J61").Value = 3081.6667
$ws.Range("K61").Value = 5047.2
$ws.Range("L61").Value = 3081.6667
$ws.Range("M61").Value = -4835.2
$ws.Range("N61").Value = -3505.6667
# Row 102
$ws.Range("H102").Value = 3455
$ws.Range("I102").Value = 3840
$ws.Range("J102").Value = 2300
$ws.Range("K102").Value = 3840
$ws.Range("L102").Value = 2300
$ws.Range("M102").Value = -2218
$ws.Range("N102").Value = -5544
# Row 110
$ws.Range("H110").Value = 1487.8572
$ws.Range("I110").Value = 1240.6666
$ws.Range("J110").Value = 1673.25
$ws.Range("K110").Value = 1240.6666
$ws.Range("L110").Value = 1673.25
$ws.Range("M110").Value = 804.3334
$ws.Range("N110").Value = -5763.25
# Row 123
$ws.Range("H123").Value = 24229.285
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 24229.285
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 24229.285
$ws.Range("N123").Value = -34029.285
# Row 132
$ws.Range("H132").Value = 3369.0735
$ws.Range("I132").Value = 3586.18
$ws.Range("J132").Value = 2766
$ws.Range("K132").Value = 10758.54
$ws.Range("L132").Value = 8298
$ws.Range("M132").Value = -8228.539999999999
$ws.Range("N132").Value = -13358
# Row 136
$ws.Range("H136").Value = 3975.0908
$ws.Range("I136").Value = 5047.2
$ws.Range("J136").Value = 3081.6667
$ws.Range("K136").Value = 15141.6
$ws.Range("L136").Value = 9245.000100000001
$ws.Range("M136").Value = -12591.6
$ws.Range("N136").Value = -14345.0001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 41
$ws.Range("H41").Value = 883.3333
$ws.Range("I41").Value = 300
$ws.Range("J41").Value = 1000
$ws.Range("K41").Value = 900
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -562
$ws.Range("N41").Value = -3676
# Row 42
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 10000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -31068
$ws.Range("M42").ClearContents()
# Row 129
$ws.Range("H129").Value = 10000873
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 10000873
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 30002619
$ws.Range("N129").Value = -30012619
$ws.Range("M129").ClearContents()
# Row 132
$ws.Range("H132").Value = 2039.4286
$ws.Range("I132").Value = 1599.4445
$ws.Range("J132").Value = 2369.4167
$ws.Range("K132").Value = 14395.0005
$ws.Range("L132").Value = 21324.7503
$ws.Range("M132").Value = -11865.0005
$ws.Range("N132").Value = -26384.7503

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2972.7273
$ws.Range("I80").Value = 2972.7273
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2972.7273
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1974.7273
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 2972.7273
$ws.Range("I83").Value = 2972.7273
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 14863.6365
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -9871.636500000001
$ws.Range("N83").ClearContents()
# Row 107
$ws.Range("H107").Value = 835.2222
$ws.Range("I107").Value = 931.5
$ws.Range("J107").Value = 498.25
$ws.Range("K107").Value = 931.5
$ws.Range("L107").Value = 498.25
$ws.Range("M107").Value = 988.5
$ws.Range("N107").Value = -4338.25
# Row 109
$ws.Range("H109").Value = 19336.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 19336.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 19336.5
$ws.Range("N109").Value = -21416.5
# Row 113
$ws.Range("H113").Value = 1872.3
$ws.Range("I113").Value = 1339.4
$ws.Range("J113").Value = 2405.2
$ws.Range("K113").Value = 1339.4
$ws.Range("L113").Value = 2405.2
$ws.Range("M113").Value = 830.5999999999999
$ws.Range("N113").Value = -6745.2
# Row 123
$ws.Range("H123").Value = 18586.908
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 18586.908
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 18586.908
$ws.Range("N123").Value = -23486.908

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 34087.332
$ws.Range("I61").Value = 34087.332
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 34087.332
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -33885.332
$ws.Range("N61").ClearContents()
# Row 64
$ws.Range("H64").Value = 31500
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 31500
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 31500
$ws.Range("N64").Value = -31950
# Row 67
$ws.Range("H67").Value = 31500
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 31500
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 31500
$ws.Range("N67").Value = -33060
# Row 68
$ws.Range("H68").Value = 3500
$ws.Range("I68").Value = 4333.3335
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 4333.3335
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -3584.3335
$ws.Range("N68").Value = -4498
# Row 71
$ws.Range("H71").Value = 3500
$ws.Range("I71").Value = 4333.3335
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 21666.6675
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -17922.6675
$ws.Range("N71").Value = -22488
# Row 76
$ws.Range("H76").Value = 20414.285
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 20414.285
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 20414.285
$ws.Range("N76").Value = -21090.285
# Row 79
$ws.Range("H79").Value = 20414.285
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 20414.285
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 20414.285
$ws.Range("N79").Value = -22754.285
# Row 113
$ws.Range("H113").Value = 34087.332
$ws.Range("I113").Value = 34087.332
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 34087.332
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -31917.332
$ws.Range("N113").ClearContents()
# Row 122
$ws.Range("H122").Value = 21431972
$ws.Range("I122").Value = 16669994
$ws.Range("J122").Value = 33336918
$ws.Range("K122").Value = 50009982
$ws.Range("L122").Value = 100010754
$ws.Range("M122").Value = -50007532
$ws.Range("N122").Value = -100015654

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1695
$ws.Range("I96").Value = 1695
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1695
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -322
$ws.Range("N96").ClearContents()
# Row 132
$ws.Range("H132").Value = 2178.375
$ws.Range("I132").Value = 1392.591
$ws.Range("J132").Value = 3138.7778
$ws.Range("K132").Value = 4177.772999999999
$ws.Range("L132").Value = 9416.3334
$ws.Range("M132").Value = -1647.772999999999
